$wb = $excel.ActiveWorkbook

# Sheet name -> last used row (in the ORIGINAL workbook, before any edits below)
$sheetLastRow = @{
    "Coiffe"    = 21
    "Cape"      = 22
    "Amulette"  = 27
    "Anneau"    = 25
    "Ceinture"  = 32
    "Bottes"    = 36
    "Arme"      = 76
    "Bouclier"  = 4
}

# 1) Sheet-specific data fix: "Anneau de la mort" row (row 13) is removed from
#    the "Anneau" sheet, shifting every following row up by one.
$wsAnneau = $wb.Worksheets.Item("Anneau")
$wsAnneau.Rows.Item(13).Delete()
$sheetLastRow["Anneau"] = 24

# Sheets whose I1 cell already carries an explicit text number format
# (numFmtId 49, "@") -- plain assignment keeps them as text. The others use
# the "General" style, so a date-like string needs a leading quote to keep
# Excel from reinterpreting it as a date serial number.
$plainTextDateSheets = @("Coiffe", "Bouclier")

foreach ($name in $sheetLastRow.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $sheetLastRow[$name]

    # 2) Header text/date update.
    $ws.Range("H1").Value = "Last save"
    if ($plainTextDateSheets -contains $name) {
        $ws.Range("I1").Value = "06-05-2021"
    } else {
        $ws.Range("I1").Value = "'06-05-2021"
    }

    # 3) "Jours consécutifs" column: drop the consecutive-day requirement,
    #    every data row goes from 1 to 0.
    $ws.Range("D2:D" + $lastRow).Value = 0

    # 4) Normalize the J:P (10-16) / Q:XFD (17-16384) column grouping: widen
    #    the first custom-width block from J:N to J:P.
    $ws.Columns.Item(14).Insert()
    $ws.Columns.Item(14).Insert()
    $ws.Columns.Item(16383).Delete()
    $ws.Columns.Item(16383).Delete()
}

# 5) Bouclier sheet selection moves from I1 to H8.
$wsBouclier = $wb.Worksheets.Item("Bouclier")
$wsBouclier.Activate()
$wsBouclier.Range("H8").Select()
